# The "seligdar.ru" entries (rows 76 and 77 - IPs 84.47.160.18 / 84.47.160.19)
# were removed from the sheet. Deleting the two entire rows shifts every
# row below them up by two and drops the now-unused shared strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("76:77").Delete()
